# Apply the calibration stage-uncertainty update:
#  - the F12 measured value (Y axis, repeatability->straightness input) was
#    re-measured: 109 -> 191 (dependent formulas in G12/H12/I12 and the
#    summary totals in G30/I30 recalculate automatically)
#  - the "Compensated" header cell (I32) gets a top border added, closing
#    off the summary box to match the neighbouring G32 cell, which
#    introduces a new border/cell style
#  - the active selection left on the sheet when the workbook was saved
#    moved to G30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data update -----------------------------------------------------
$ws.Range("F12").Value = 191

# --- Formatting update -------------------------------------------------
# Add a thin top border to I32 (xlEdgeTop = 8), matching the new border
# style added to the workbook (medium/medium/thin/medium).
$iEdgeTop = 8
$iThin = 2
$iContinuous = 1
$ws.Range("I32").Borders.Item($iEdgeTop).LineStyle = $iContinuous
$ws.Range("I32").Borders.Item($iEdgeTop).Weight = $iThin

# --- Selection update ---------------------------------------------------
$ws.Range("G30").Select()
